# level 7 stuff, change house icons, added ranking
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for "bonus" and "total" key/value pairs, just above
# the "cube_field" row (which is row 27 in the original layout).
$ws.Rows("27:28").Insert()

$ws.Range("A27").Value = "bonus"
$ws.Range("B27").Value = "Bonus:"
$ws.Range("A28").Value = "total"
$ws.Range("B28").Value = "Total:"

# Update the house / field / pond display names (rows shifted down by 2
# after the insert above). Set in this particular order so new shared
# strings are interned in the same order as the source workbook.
$ws.Range("B33").Value = "Pennyroyal House"
$ws.Range("B31").Value = "Marigold House"
$ws.Range("B29").Value = "Pleasant Field"
$ws.Range("B30").Value = "Serene Pond"
$ws.Range("B32").Value = "Green House"

# Update the view so it matches the author's final selection state.
$ws.Range("B32").Select()
